$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Secondary Income (year) - replace hard-coded 0 with a formula
$ws.Range("B3").Formula = "=23.8*52*10"

# Student Loan Repayment (Month) - guard the repayment formula so it
# is zero when there's no student loan balance left (B4)
$ws.Range("B17").Formula = "=IF(B4>0,0.12*(`$B`$2/12-1690) + 0.12*B3/26*4,0)"

# Move the active selection to B5 (reflects where the user left off editing)
$ws.Range("B5").Select()

$wb.Save()
